$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Adicionada REC de ALP" - fill in Recuperacao 02 (column F) scores
# for the students that took the make-up exam.
$ws.Range("F12").Value = 10
$ws.Range("F15").Value = 4
$ws.Range("F17").Value = 8
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = 8
$ws.Range("F23").Value = 8
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = 8
$ws.Range("F28").Value = 6
$ws.Range("F30").Value = 8
$ws.Range("F35").Value = 6
$ws.Range("F36").Value = 6

$ws.Range("F12").Select()
